$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 28
$ws.Cells.Item(28, 4).Value = 44463   # D28: 44166 -> 44463
$ws.Cells.Item(28, 13).Value = 240   # M28: 300 -> 240
$ws.Cells.Item(28, 14).Value = 2600   # N28: 14000 -> 2600
$ws.Cells.Item(28, 15).Value = 2700   # O28: 14500 -> 2700
$ws.Cells.Item(28, 16).Value = 2650   # P28: 14250 -> 2650
$ws.Cells.Item(28, 17).Value = '$/kilo (en caja de 15 kilos)'   # Q28: '$/bandeja 8 kilos' -> '$/kilo (en caja de 15 kilos)'
$ws.Cells.Item(28, 19).Value = 2650   # S28: 1781 -> 2650
$ws.Cells.Item(28, 20).Value = 1   # T28: 8 -> 1

# Row 29
$ws.Cells.Item(29, 4).Value = 44463   # D29: 44166 -> 44463
$ws.Cells.Item(29, 13).Value = 300   # M29: 200 -> 300
$ws.Cells.Item(29, 14).Value = 2200   # N29: 12000 -> 2200
$ws.Cells.Item(29, 15).Value = 2300   # O29: 12500 -> 2300
$ws.Cells.Item(29, 16).Value = 2250   # P29: 12250 -> 2250
$ws.Cells.Item(29, 17).Value = '$/kilo (en caja de 15 kilos)'   # Q29: '$/bandeja 8 kilos' -> '$/kilo (en caja de 15 kilos)'
$ws.Cells.Item(29, 19).Value = 2250   # S29: 1531 -> 2250
$ws.Cells.Item(29, 20).Value = 1   # T29: 8 -> 1

# Row 30
$ws.Cells.Item(30, 4).Value = 44463   # D30: 44168 -> 44463
$ws.Cells.Item(30, 12).Value = 'Segunda'   # L30: 'Especial' -> 'Segunda'
$ws.Cells.Item(30, 14).Value = 1900   # N30: 14000 -> 1900
$ws.Cells.Item(30, 15).Value = 2000   # O30: 14500 -> 2000
$ws.Cells.Item(30, 16).Value = 1950   # P30: 14250 -> 1950
$ws.Cells.Item(30, 17).Value = '$/kilo (en caja de 15 kilos)'   # Q30: '$/bandeja 8 kilos' -> '$/kilo (en caja de 15 kilos)'
$ws.Cells.Item(30, 19).Value = 1950   # S30: 1781 -> 1950
$ws.Cells.Item(30, 20).Value = 1   # T30: 8 -> 1

# Row 31
$ws.Cells.Item(31, 4).Value = 44166   # D31: 44168 -> 44166
$ws.Cells.Item(31, 12).Value = 'Especial'   # L31: 'Primera' -> 'Especial'
$ws.Cells.Item(31, 13).Value = 300   # M31: 200 -> 300
$ws.Cells.Item(31, 14).Value = 14000   # N31: 12000 -> 14000
$ws.Cells.Item(31, 15).Value = 14500   # O31: 12500 -> 14500
$ws.Cells.Item(31, 16).Value = 14250   # P31: 12250 -> 14250
$ws.Cells.Item(31, 19).Value = 1781   # S31: 1531 -> 1781

# Row 32
$ws.Cells.Item(32, 4).Value = 44166   # D32: 44168 -> 44166
$ws.Cells.Item(32, 12).Value = 'Primera'   # L32: 'Segunda' -> 'Primera'
$ws.Cells.Item(32, 14).Value = 12000   # N32: 9500 -> 12000
$ws.Cells.Item(32, 15).Value = 12500   # O32: 10000 -> 12500
$ws.Cells.Item(32, 16).Value = 12250   # P32: 9750 -> 12250
$ws.Cells.Item(32, 19).Value = 1531   # S32: 1219 -> 1531

# Row 33
$ws.Cells.Item(33, 4).Value = 44168   # D33: 44162 -> 44168
$ws.Cells.Item(33, 13).Value = 240   # M33: 340 -> 240

# Row 34
$ws.Cells.Item(34, 4).Value = 44168   # D34: 44162 -> 44168
$ws.Cells.Item(34, 13).Value = 200   # M34: 300 -> 200

# Row 35
$ws.Cells.Item(35, 4).Value = 44168   # D35: 44162 -> 44168

# Row 36
$ws.Cells.Item(36, 4).Value = 44162   # D36: 44410 -> 44162
$ws.Cells.Item(36, 12).Value = 'Especial'   # L36: 'Primera' -> 'Especial'
$ws.Cells.Item(36, 13).Value = 340   # M36: 240 -> 340
$ws.Cells.Item(36, 14).Value = 14000   # N36: 2400 -> 14000
$ws.Cells.Item(36, 15).Value = 14500   # O36: 2500 -> 14500
$ws.Cells.Item(36, 16).Value = 14250   # P36: 2450 -> 14250
$ws.Cells.Item(36, 17).Value = '$/bandeja 8 kilos'   # Q36: '$/kilo (en caja de 15 kilos)' -> '$/bandeja 8 kilos'
$ws.Cells.Item(36, 18).Value = 'Provincia de Limarí'   # R36: 'Provincia del Elquí' -> 'Provincia de Limarí'
$ws.Cells.Item(36, 19).Value = 1781   # S36: 2450 -> 1781
$ws.Cells.Item(36, 20).Value = 8   # T36: 1 -> 8

# Row 37
$ws.Cells.Item(37, 4).Value = 44162   # D37: 44410 -> 44162
$ws.Cells.Item(37, 12).Value = 'Primera'   # L37: 'Segunda' -> 'Primera'
$ws.Cells.Item(37, 13).Value = 300   # M37: 240 -> 300
$ws.Cells.Item(37, 14).Value = 12000   # N37: 2000 -> 12000
$ws.Cells.Item(37, 15).Value = 12500   # O37: 2100 -> 12500
$ws.Cells.Item(37, 16).Value = 12250   # P37: 2050 -> 12250
$ws.Cells.Item(37, 17).Value = '$/bandeja 8 kilos'   # Q37: '$/kilo (en caja de 15 kilos)' -> '$/bandeja 8 kilos'
$ws.Cells.Item(37, 18).Value = 'Provincia de Limarí'   # R37: 'Provincia del Elquí' -> 'Provincia de Limarí'
$ws.Cells.Item(37, 19).Value = 1531   # S37: 2050 -> 1531
$ws.Cells.Item(37, 20).Value = 8   # T37: 1 -> 8

# Row 38
$ws.Cells.Item(38, 4).Value = 44162   # D38: 44410 -> 44162
$ws.Cells.Item(38, 12).Value = 'Segunda'   # L38: 'Tercera' -> 'Segunda'
$ws.Cells.Item(38, 14).Value = 9500   # N38: 1600 -> 9500
$ws.Cells.Item(38, 15).Value = 10000   # O38: 1700 -> 10000
$ws.Cells.Item(38, 16).Value = 9750   # P38: 1650 -> 9750
$ws.Cells.Item(38, 17).Value = '$/bandeja 8 kilos'   # Q38: '$/kilo (en caja de 15 kilos)' -> '$/bandeja 8 kilos'
$ws.Cells.Item(38, 18).Value = 'Provincia de Limarí'   # R38: 'Provincia del Elquí' -> 'Provincia de Limarí'
$ws.Cells.Item(38, 19).Value = 1219   # S38: 1650 -> 1219
$ws.Cells.Item(38, 20).Value = 8   # T38: 1 -> 8

# Row 39
$ws.Cells.Item(39, 4).Value = 44410   # D39: 44411 -> 44410
$ws.Cells.Item(39, 13).Value = 240   # M39: 600 -> 240

# Row 40
$ws.Cells.Item(40, 4).Value = 44410   # D40: 44411 -> 44410
$ws.Cells.Item(40, 13).Value = 240   # M40: 400 -> 240

# Row 41 (new)
$ws.Cells.Item(41, 1).Value = 8   # A41
$ws.Cells.Item(41, 2).Value = 'Terminal La Palmera de La Serena'   # B41
$ws.Cells.Item(41, 3).Value = 'Coquimbo'   # C41
$ws.Cells.Item(41, 4).Value = 44410   # D41
$ws.Cells.Item(41, 5).Value = 4   # E41
$ws.Cells.Item(41, 6).Value = 'Fruta'   # F41
$ws.Cells.Item(41, 7).Value = 100107   # G41
$ws.Cells.Item(41, 8).Value = 'Otros'   # H41
$ws.Cells.Item(41, 9).Value = 100107002   # I41
$ws.Cells.Item(41, 10).Value = 'Chirimoya'   # J41
$ws.Cells.Item(41, 11).Value = 'Cultivar IV Región'   # K41
$ws.Cells.Item(41, 12).Value = 'Tercera'   # L41
$ws.Cells.Item(41, 13).Value = 200   # M41
$ws.Cells.Item(41, 14).Value = 1600   # N41
$ws.Cells.Item(41, 15).Value = 1700   # O41
$ws.Cells.Item(41, 16).Value = 1650   # P41
$ws.Cells.Item(41, 17).Value = '$/kilo (en caja de 15 kilos)'   # Q41
$ws.Cells.Item(41, 18).Value = 'Provincia del Elquí'   # R41
$ws.Cells.Item(41, 19).Value = 1650   # S41
$ws.Cells.Item(41, 20).Value = 1   # T41
$ws.Cells.Item(41, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"   # D41 date style

# Row 42 (new)
$ws.Cells.Item(42, 1).Value = 8   # A42
$ws.Cells.Item(42, 2).Value = 'Terminal La Palmera de La Serena'   # B42
$ws.Cells.Item(42, 3).Value = 'Coquimbo'   # C42
$ws.Cells.Item(42, 4).Value = 44411   # D42
$ws.Cells.Item(42, 5).Value = 4   # E42
$ws.Cells.Item(42, 6).Value = 'Fruta'   # F42
$ws.Cells.Item(42, 7).Value = 100107   # G42
$ws.Cells.Item(42, 8).Value = 'Otros'   # H42
$ws.Cells.Item(42, 9).Value = 100107002   # I42
$ws.Cells.Item(42, 10).Value = 'Chirimoya'   # J42
$ws.Cells.Item(42, 11).Value = 'Cultivar IV Región'   # K42
$ws.Cells.Item(42, 12).Value = 'Primera'   # L42
$ws.Cells.Item(42, 13).Value = 600   # M42
$ws.Cells.Item(42, 14).Value = 2400   # N42
$ws.Cells.Item(42, 15).Value = 2500   # O42
$ws.Cells.Item(42, 16).Value = 2450   # P42
$ws.Cells.Item(42, 17).Value = '$/kilo (en caja de 15 kilos)'   # Q42
$ws.Cells.Item(42, 18).Value = 'Provincia del Elquí'   # R42
$ws.Cells.Item(42, 19).Value = 2450   # S42
$ws.Cells.Item(42, 20).Value = 1   # T42
$ws.Cells.Item(42, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"   # D42 date style

# Row 43 (new)
$ws.Cells.Item(43, 1).Value = 8   # A43
$ws.Cells.Item(43, 2).Value = 'Terminal La Palmera de La Serena'   # B43
$ws.Cells.Item(43, 3).Value = 'Coquimbo'   # C43
$ws.Cells.Item(43, 4).Value = 44411   # D43
$ws.Cells.Item(43, 5).Value = 4   # E43
$ws.Cells.Item(43, 6).Value = 'Fruta'   # F43
$ws.Cells.Item(43, 7).Value = 100107   # G43
$ws.Cells.Item(43, 8).Value = 'Otros'   # H43
$ws.Cells.Item(43, 9).Value = 100107002   # I43
$ws.Cells.Item(43, 10).Value = 'Chirimoya'   # J43
$ws.Cells.Item(43, 11).Value = 'Cultivar IV Región'   # K43
$ws.Cells.Item(43, 12).Value = 'Segunda'   # L43
$ws.Cells.Item(43, 13).Value = 400   # M43
$ws.Cells.Item(43, 14).Value = 2000   # N43
$ws.Cells.Item(43, 15).Value = 2100   # O43
$ws.Cells.Item(43, 16).Value = 2050   # P43
$ws.Cells.Item(43, 17).Value = '$/kilo (en caja de 15 kilos)'   # Q43
$ws.Cells.Item(43, 18).Value = 'Provincia del Elquí'   # R43
$ws.Cells.Item(43, 19).Value = 2050   # S43
$ws.Cells.Item(43, 20).Value = 1   # T43
$ws.Cells.Item(43, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"   # D43 date style
